$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find a paragraph whose exact text (minus trailing paragraph mark)
# matches the given string.
# ---------------------------------------------------------------------------
function Find-ParagraphByExactText($doc, [string]$text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t -eq ($text + "`r")) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Merge the "Notice u/s 94 BNSS, 2023" paragraph into the "To," paragraph
#    so both sit on one line, separated by a centered tab stop at 4680 twips
#    (234 pt), with "Notice u/s 94 BNSS, 2023" rendered bold.
# ---------------------------------------------------------------------------
$noticeText = "Notice u/s 94 BNSS, 2023"
$toText = "To,"

$pNotice = Find-ParagraphByExactText $d $noticeText
$pTo = Find-ParagraphByExactText $d $toText

# Delete the whole "Notice u/s 94 BNSS, 2023" paragraph, including its
# paragraph mark. This merges it away; the surviving paragraph keeps the
# "To," paragraph's own pPr/formatting.
$fullNotice = $d.Range($pNotice.Range.Start, $pNotice.Range.End)
$fullNotice.Delete()

# Re-acquire the (now shifted) "To," paragraph.
$pTo = Find-ParagraphByExactText $d $toText

# Add a centered custom tab stop at 4680 twips = 234 points.
$tabs = $pTo.Format.TabStops
$tabs.Add(234, 1)

# Clear the existing run content of the "To," paragraph (but keep the
# paragraph mark itself), so re-typed text starts with clean formatting.
$bodyRange = $d.Range($pTo.Range.Start, $pTo.Range.End - 1)
if ($bodyRange.Start -ne $bodyRange.End) {
    $bodyRange.Delete()
}

# Re-insert "To," with no direct character formatting.
$insertPt = $d.Range($pTo.Range.Start, $pTo.Range.Start)
$insertPt.InsertAfter($toText)

# Insert a tab right after "To,".
$pTo = $d.Paragraphs.Item($pTo.Index)
$tabPt = $d.Range($pTo.Range.End - 1, $pTo.Range.End - 1)
$tabPt.InsertAfter([char]9)

# Insert "Notice u/s 94 BNSS, 2023" after the tab.
$pTo = $d.Paragraphs.Item($pTo.Index)
$noticeStart = $pTo.Range.End - 1
$noticeInsertPt = $d.Range($noticeStart, $noticeStart)
$noticeInsertPt.InsertAfter($noticeText)

# Bold only the newly-inserted "Notice u/s 94 BNSS, 2023" text.
$pTo = $d.Paragraphs.Item($pTo.Index)
$noticeRange = $d.Range($noticeStart, $pTo.Range.End - 1)
$noticeRange.Bold = 1

# ---------------------------------------------------------------------------
# 2) Add w:before="0" (SpaceBefore = 0) to the four numbered request
#    paragraphs.
# ---------------------------------------------------------------------------
$targets = @(
    "1. Details of the user (Name, Address, Contact No. etc.) to whom below IP’s were allotted at the mentioned Date & time against each.",
    "2. Kindly provide the ownership of the users, to whom IP was allotted.",
    "3. Kindly preserve the record till further directions.  ",
    "4. Kindly provide any other useful details. "
)

foreach ($target in $targets) {
    $p = Find-ParagraphByExactText $d $target
    if ($p -ne $null) {
        $p.Format.SpaceBefore = 0
    }
}
